$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the "How are you
# doing." paragraph. Remove it -- it will be re-created at the end of the
# newly typed paragraph below, matching Word's normal behaviour of moving
# _GoBack to the most recent edit point.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Append a new paragraph after the last one ("How are you doing.") and
# type the new sentence into it.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Im fine thanks"

# Re-create the _GoBack bookmark at the end of the new paragraph (i.e. at
# the end of the text we just typed), mirroring what Word does whenever
# new text is entered.
$newRange = $newPara.Range
$newRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $newRange)
